# The source data range ("NEW" sheet) had a duplicated record at row 16
# (Caso 5460 / MOLDES 2735) that needs to be removed. Deleting the entire
# row shifts every subsequent record up by one, which matches the newer
# export (last row, old 57 / Caso 6392, disappears and the sheet's used
# range shrinks from A1:P57 to A1:P56).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(16).Delete()
